# Populate row 11 of the Kayserispor match-log sheet with the
# 24/10/2025 Karagumruk vs Kayserispor result and its associated stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "24/10/2025"
$ws.Range("B11").Value = "Karagumruk"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = "Kayserispor"
$ws.Range("F11").Value = "D"
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 2
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 1.66
$ws.Range("L11").Value = 2.04
$ws.Range("M11").Value = 14
$ws.Range("N11").Value = 10
$ws.Range("O11").Value = 10
$ws.Range("P11").Value = 3
